$d = $word.ActiveDocument

# The document's header/footer stories each carry a single inline picture
# (the Pearson logo in the footers, the BTEC logo in the headers). Renaming
# a picture in Word updates the drawing's display name (wp:docPr/@name,
# surfaced on the object model as InlineShape.Name) while leaving its
# alt-text (descr) and numeric id untouched.
#
#   Footers: PearsonLogo picture "image1.png" -> "image2.png"
#   Headers: BTec_Logo-Orange picture "image2.jpg" -> "image1.jpg"

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)

    for ($i = 1; $i -le 3; $i++) {
        $hdr = $section.Headers.Item($i)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($k = 1; $k -le $shapes.Count; $k++) {
                $shape = $shapes.Item($k)
                if ($shape.AlternativeText -eq "BTec_Logo-Orange") {
                    $shape.Name = "image1.jpg"
                }
            }
        }

        $ftr = $section.Footers.Item($i)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($k = 1; $k -le $shapes.Count; $k++) {
                $shape = $shapes.Item($k)
                if ($shape.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shape.Name = "image2.png"
                }
            }
        }
    }
}
